$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, matching the formatting of the
# existing header cells (bold/centered/bordered style used by B1:G1).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the corresponding data value for row 2 (plain numeric cell, same as
# the other data cells).
$ws.Range("H2").Value = 0
